$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in column B (risk description) first, top to bottom,
#     so the shared-string table gets these entries in row order ---
$ws.Range("B2").Value = "Funcionário faltar"
$ws.Range("B3").Value = "problemas na infraestrutura "
$ws.Range("B4").Value = "Perda de informações e dados essenciais "
$ws.Range("B5").Value = "Ajuestes não realizados dentro do prazo"
$ws.Range("B6").Value = "Reclamações do usuário devido a bugs"

# --- Probabilidade (C) / Impacto (D) ---
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 1

$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 3

$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 3

$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 3

$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 2

# --- Ação (F) ---
$ws.Range("F2").Value = "Mitigar"
$ws.Range("F3").Value = "Evitar"
$ws.Range("F4").Value = "Evitar"
$ws.Range("F5").Value = "Evitar"
$ws.Range("F6").Value = "Evitar"

# --- Como? (G), filled after column B/F entries so the new shared
#     strings are appended after them, matching the authoring order ---
$ws.Range("G2").Value = "Nova adequação e redistribuição da equipe"
$ws.Range("G3").Value = "Fazer analise preventiva todo dia "
$ws.Range("G4").Value = "Ter protocolo de backups após todas atualizações"
$ws.Range("G5").Value = "Ter cronograma e planejamento de odo o projeto"
$ws.Range("G6").Value = "Garantir que a homologação foi feita de maneira acertiva]"

# --- Match column G formatting (center / center) with the rest of the
#     row, reusing F's existing style (copy format only, so no new
#     style entry is created in styles.xml) ---
$ws.Range("F2").Copy()
$ws.Range("G2:G11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row heights for the filled rows ---
$ws.Rows.Item(2).RowHeight = 29.25
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 22.5
$ws.Rows.Item(6).RowHeight = 27

# --- Column G width widened to 53 ---
$ws.Columns.Item(7).ColumnWidth = 52.14

# --- G13: empty cell, underline font applied (no alignment change) ---
$ws.Range("G13").Font.Underline = $true

# --- Sheet view: selection changes (final active cell G13) ---
$ws.Range("G13").Select()
